$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Feria Lagunitas de Puerto Montt - Cebollín) was
# recorded. Insert it as row 226, pushing the existing rows 226..259 down to
# 227..260 (mirrors Excel's "insert row" behaviour, carrying formatting such
# as the date style in column D along with it).
$ws.Rows(226).Insert()

$ws.Cells.Item(226, 1).Value = 4
$ws.Cells.Item(226, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(226, 3).Value = "Los Lagos"
$ws.Cells.Item(226, 4).Value = 44694
$ws.Cells.Item(226, 5).Value = 10
$ws.Cells.Item(226, 6).Value = 100112037
$ws.Cells.Item(226, 7).Value = "Cebollín"
$ws.Cells.Item(226, 8).Value = "Sin especificar"
$ws.Cells.Item(226, 9).Value = "Primera"
$ws.Cells.Item(226, 10).Value = 120
$ws.Cells.Item(226, 11).Value = 11000
$ws.Cells.Item(226, 12).Value = 11000
$ws.Cells.Item(226, 13).Value = 11000
$ws.Cells.Item(226, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(226, 15).Value = "Región Metropolitana"
$ws.Cells.Item(226, 16).Value = 306
$ws.Cells.Item(226, 17).Value = 36
$ws.Cells.Item(226, 18).Value = "Hortaliza"
